$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.638.91'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '1.849.29'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.31'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4234'
$ws.Range("E7").Value = '  -1.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3641'
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.40'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07295'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8759'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.69'
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").Value = '1.833.00'
$ws.Range("E13").Value = '  -0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.339'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.524'
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06869'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.62'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008896'
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.36'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = '27.670.02'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("E23").Value = '  -1.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.38'
$ws.Range("E24").Value = '  -5.21%  '
$ws.Range("D25").Value = '2.100.09'
$ws.Range("E25").Value = '  -2.90%  '
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.30'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.95'
$ws.Range("E28").Value = '  +3.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '122.01'
$ws.Range("E29").Value = '  +8.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.254'
$ws.Range("E30").Value = '  -3.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.878'
$ws.Range("E31").Value = '  +12.16%  '
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7674'
$ws.Range("E33").Value = '  -3.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.549'
$ws.Range("E34").Value = '  -3.56%  '
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.107'
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05346'
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01930'
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.813'
$ws.Range("E41").Value = '  -5.74%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5107'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.887'
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1647'
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.308'
$ws.Range("E45").Value = '  -4.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06536'
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.58'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.31'
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4683'
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.627'
$ws.Range("E51").Value = '  -1.79%  '
